$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that was on A72 (http://epp.eurostat... link)
$ws.Range("A72").Hyperlinks.Delete()

# --- Reorganize the "Source" block (rows 70-74) ---
# Row 70 "Source:" stays as-is.
# Row 71 used to hold the SBS Main Indicators description; it becomes blank.
$ws.Range("A71").Value = ""

# Row 72 used to hold the hyperlink URL; it now holds the SBS Main Indicators text (no hyperlink, plain "source" style).
$ws.Range("A72").Value = "SBS Main Indicators, Annual enterprise statistics by size class for special aggregates of activities (NACE Rev. 2)"
$ws.Range("A72").Style = "source"

# Row 73 (empty) stays as-is (untouched).

# Row 74 (new) now holds the URL text as plain text, "source" style.
$ws.Range("A74").Value = "http://epp.eurostat.ec.europa.eu/portal/page/portal/european_business/data/database"
$ws.Range("A74").Style = "source"

# --- Reorganize the I.Stat / SBS Eurostat block ---
# Row 76 "I.Stat (Italian Official Statistics)" (title) moves down to row 77.
$ws.Range("A77").Value = "I.Stat (Italian Official Statistics)"
$ws.Range("A77").Style = "title"

# Row 77's long description is replaced by a duplicate of the I.Stat title text, "source" style, at row 78.
$ws.Range("A78").Value = "I.Stat (Italian Official Statistics)"
$ws.Range("A78").Style = "source"

# Row 78 "SBS Eurostat" (title) moves down to row 79.
$ws.Range("A79").Value = "SBS Eurostat"
$ws.Range("A79").Style = "title"

# Row 79's long description is replaced by a duplicate of "SBS Eurostat", "source" style, at row 80.
$ws.Range("A80").Value = "SBS Eurostat"
$ws.Range("A80").Style = "source"

# Clear row 76 (now unused, since its content moved to row 77)
$ws.Range("A76").Clear()
